# Update countries & provincias Spain
#
# The source data table (sheet "Pais") is refreshed with newer COVID-19
# figures. Because the table is kept sorted by total cases (column B,
# descending), two blocks of rows need their country label and figures
# re-assigned to reflect the new ranking, and the "last updated" timestamp
# in A1 is bumped.
#
# Block 1 (rows 7-8): "Estados Unidos" overtakes "Alemania".
# Block 2 (rows 87-91): "Lituania" jumps to the top of this block, pushing
#   "Bielorrusia", "Moldavia", "Venezuela" and "Malta" down one row each.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CountryRow($row, $country, $values) {
    $ws.Cells.Item($row, 1).Value = $country
    for ($i = 0; $i -lt $values.Length; $i++) {
        $ws.Cells.Item($row, 2 + $i).Value = $values[$i]
    }
}

# --- "Datos actualizados..." timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 21 de Marzo de 2020 a las 06:16"

# --- Rows 7-8 ---
Set-CountryRow 7 "Estados Unidos" @(20193, 810, 147, 19766, 64, 24, 280)
Set-CountryRow 8 "Alemania"       @(19848, 0,   180, 19600, 2,  0,  68)

# --- Rows 87-91 ---
Set-CountryRow 87 "Lituania"    @(69, 0, 1,  67, 1, 0, 1)
Set-CountryRow 88 "Bielorrusia" @(69, 0, 15, 54, 0, 0, 0)
Set-CountryRow 89 "Moldavia"    @(66, 0, 1,  64, 3, 0, 1)
Set-CountryRow 90 "Venezuela"   @(65, 0, 1,  64, 0, 0, 0)
Set-CountryRow 91 "Malta"       @(64, 0, 2,  62, 1, 0, 0)
